# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to Sheets workbook (Aegis_Profits)
$wb = $excel.ActiveWorkbook

$updates = @{}

$updates["ALC"] = @(
    @("H15", 1193.2439),
    @("I15", 1193.2439),
    @("K15", 3579.7317),
    @("M15", -3410.7317),
    @("H17", 967.7353000000001),
    @("J17", 967.7353000000001),
    @("L17", 2903.2059),
    @("N17", -3239.2059),
    @("H53", 546.43475),
    @("I53", 678.6667),
    @("K53", 678.6667),
    @("M53", -41.66669999999999),
    @("H112", 1291.4117),
    @("I112", 941.4286),
    @("J112", 1536.4),
    @("K112", 2824.2858),
    @("L112", 4609.200000000001),
    @("M112", -1716.2858),
    @("N112", -6825.200000000001),
    @("H125", 1403.1666),
    @("I125", 1192.4286),
    @("J125", 1587.5625),
    @("K125", 10731.8574),
    @("L125", 14288.0625),
    @("M125", -8271.857399999999),
    @("N125", -19208.0625),
    @("H129", 2571.7542),
    @("J129", 1200.2683),
    @("L129", 3600.8049),
    @("N129", -13600.8049),
    @("H132", 10878695),
    @("I132", 14717904),
    @("J132", 934),
    @("K132", 44153712),
    @("L132", 2802),
    @("M132", -44151182),
    @("N132", -7862),
    @("H137", 990.4068),
    @("I137", 971.9636),
    @("K137", 2915.8908),
    @("M137", -365.8908000000001),
    @("H138", 1853.9166),
    @("I138", 1390.8572),
    @("J138", 3474.625),
    @("K138", 4172.571599999999),
    @("L138", 10423.875),
    @("M138", 967.4284000000007),
    @("N138", -20703.875),
    @("H141", 1327.8422),
    @("I141", 1267.625),
    @("K141", 3802.875),
    @("M141", 1377.125)
)

$updates["ARM"] = @(
    @("H4", 193),
    @("I4", 100),
    @("J4", 239.5),
    @("K4", 100),
    @("L4", 239.5),
    @("M4", 16),
    @("N4", -471.5),
    @("H11", 49800),
    @("J11", 49800),
    @("L11", 49800),
    @("N11", -50088),
    @("H32", 20516.553),
    @("I32", 4533.036),
    @("J32", 141119.45),
    @("K32", 4533.036),
    @("L32", 141119.45),
    @("M32", -4246.036),
    @("N32", -141693.45),
    @("H61", 1361.1945),
    @("I61", 1114.875),
    @("K61", 1114.875),
    @("M61", -902.875),
    @("H74", 584.2941),
    @("I74", 542.86664),
    @("J74", 895),
    @("K74", 542.86664),
    @("L74", 895),
    @("M74", 331.13336),
    @("N74", -2643),
    @("H77", 584.2941),
    @("I77", 542.86664),
    @("J77", 895),
    @("K77", 2714.3332),
    @("L77", 4475),
    @("M77", 1653.6668),
    @("N77", -13211),
    @("H132", 4633.5557),
    @("I132", 4755.5864),
    @("J132", 4128),
    @("K132", 14266.7592),
    @("L132", 12384),
    @("M132", -11736.7592),
    @("N132", -17444),
    @("H136", 1361.1945),
    @("I136", 1114.875),
    @("K136", 3344.625),
    @("M136", -794.625)
)

$updates["BSM"] = @(
    @("H134", 2411.2449),
    @("I134", 2120.6738),
    @("J134", 6866.6665),
    @("K134", 6362.0214),
    @("L134", 20599.9995),
    @("M134", -3827.0214),
    @("N134", -25669.9995)
)

$updates["CRP"] = @(
    @("H16", 1305.8182),
    @("I16", 859.375),
    @("J16", 2496.3333),
    @("K16", 859.375),
    @("L16", 2496.3333),
    @("M16", -572.375),
    @("N16", -3070.3333),
    @("H31", 24990.484),
    @("I31", 1427.8572),
    @("J31", 55534.63),
    @("K31", 1427.8572),
    @("L31", 55534.63),
    @("M31", -1132.8572),
    @("N31", -56124.63),
    @("H34", 24990.484),
    @("I34", 1427.8572),
    @("J34", 55534.63),
    @("K34", 1427.8572),
    @("L34", 55534.63),
    @("M34", -1225.8572),
    @("N34", -55938.63),
    @("H113", 1305.8182),
    @("I113", 859.375),
    @("J113", 2496.3333),
    @("K113", 859.375),
    @("L113", 2496.3333),
    @("M113", 1310.625),
    @("N113", -6836.3333),
    @("H134", 1321.591),
    @("I134", 1194.6428),
    @("J134", 1543.75),
    @("K134", 3583.9284),
    @("L134", 4631.25),
    @("M134", -1048.9284),
    @("N134", -9701.25),
    @("H141", 149799.8),
    @("J141", 149499.5),
    @("L141", 149499.5),
    @("N141", -159859.5)
)

$updates["CUL"] = @(
    @("H5", 960.75),
    @("I5", 648.55554),
    @("K5", 1945.66662),
    @("M5", -1833.66662),
    @("H17", 1000),
    @("I17", 1000),
    @("K17", 3000),
    @("M17", -2831),
    @("H100", 2300),
    @("J100", 2300),
    @("L100", 6900),
    @("N100", -8522),
    @("H122", 509.54544),
    @("I122", 494.375),
    @("J122", 550),
    @("K122", 4449.375),
    @("L122", 4950),
    @("M122", -1999.375),
    @("N122", -9850),
    @("H131", 8540),
    @("J131", 8586.579),
    @("L131", 25759.737),
    @("N131", -35839.737),
    @("H135", 960.75),
    @("I135", 648.55554),
    @("K135", 5836.99986),
    @("M135", -3301.99986)
)

$updates["GSM"] = @(
    @("H113", 1657.3334),
    @("J113", 1528.5385),
    @("L113", 1528.5385),
    @("N113", -5868.538500000001)
)

$updates["LTW"] = @(
    @("H6", 24095),
    @("I6", 18000),
    @("J6", 24965.715),
    @("K6", 18000),
    @("L6", 24965.715),
    @("M6", -17888),
    @("N6", -25189.715),
    @("H16", 3074018.5),
    @("I16", 3706531.2),
    @("J16", 1813.4286),
    @("K16", 3706531.2),
    @("L16", 1813.4286),
    @("M16", -3706361.2),
    @("N16", -2153.4286),
    @("H46", 844049.7),
    @("I46", 429.33334),
    @("J46", 1125256.5),
    @("K46", 429.33334),
    @("L46", 1125256.5),
    @("M46", -241.33334),
    @("N46", -1125632.5)
)

$updates["WVR"] = @(
    @("H47", 8000),
    @("J47", 8000),
    @("L47", 8000),
    @("N47", -9144),
    @("H126", 1349.2084),
    @("I126", 1297.8235),
    @("J126", 1474),
    @("K126", 3893.4705),
    @("L126", 4422),
    @("M126", -1423.4705),
    @("N126", -9362),
    @("H132", 1765),
    @("I132", 1670.75),
    @("K132", 5012.25),
    @("M132", -2482.25),
    @("H136", 542.5345),
    @("I136", 371.1915),
    @("K136", 1113.5745),
    @("M136", 1436.4255)
)

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($pair in $updates[$sheetName]) {
        $cellRef = $pair[0]
        $newValue = $pair[1]
        $ws.Range($cellRef).Value = $newValue
    }
}
